$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.028415333333333
$ws.Range("H2").Value = 3.085246
$ws.Range("I2").Value = 0.04565156193945813
$ws.Range("J2").Value = 0.04565156193945813
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 142.563596
$ws.Range("N2").Value = 427.690788
$ws.Range("O2").Value = 0.9315318106005318
$ws.Range("P2").Value = 0.9315318106005319
$ws.Range("Q2").Value = 146.6145881015387
$ws.Range("R2").Value = 1319.531292913848
$ws.Range("S2").Value = 0.04252588215020576
$ws.Range("T2").Value = 0.04252588215020576

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.028415333333333
$ws.Range("H3").Value = 3.085246
$ws.Range("I3").Value = 0.04565156193945813
$ws.Range("J3").Value = 0.04565156193945813
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.4502283333333333
$ws.Range("N3").Value = 1.350685
$ws.Range("O3").Value = 0.002941859116219682
$ws.Range("P3").Value = 0.002941859116219682
$ws.Range("Q3").Value = 0.4630217215011111
$ws.Range("R3").Value = 4.16719549351
$ws.Range("S3").Value = 0.0001343004636612624
$ws.Range("T3").Value = 0.0001343004636612623

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.028415333333333
$ws.Range("H4").Value = 3.085246
$ws.Range("I4").Value = 0.04565156193945813
$ws.Range("J4").Value = 0.04565156193945813
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.398212333333333
$ws.Range("N4").Value = 4.194637
$ws.Range("O4").Value = 0.009136128037020016
$ws.Range("P4").Value = 0.009136128037020014
$ws.Range("Q4").Value = 1.437943002855778
$ws.Range("R4").Value = 12.941487025702
$ws.Range("S4").Value = 0.0004170785149688393
$ws.Range("T4").Value = 0.0004170785149688392

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.028415333333333
$ws.Range("H5").Value = 3.085246
$ws.Range("I5").Value = 0.04565156193945813
$ws.Range("J5").Value = 0.04565156193945813
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.666884333333333
$ws.Range("N5").Value = 8.000653
$ws.Range("O5").Value = 0.01742582020512581
$ws.Range("P5").Value = 0.0174258202051258
$ws.Range("Q5").Value = 2.742664740626445
$ws.Range("R5").Value = 24.683982665638
$ws.Range("S5").Value = 0.0007955159104401617
$ws.Range("T5").Value = 0.0007955159104401616

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.028415333333333
$ws.Range("H6").Value = 3.085246
$ws.Range("I6").Value = 0.04565156193945813
$ws.Range("J6").Value = 0.04565156193945813
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.688331
$ws.Range("N6").Value = 11.064993
$ws.Range("O6").Value = 0.02410010515253888
$ws.Range("P6").Value = 0.02410010515253888
$ws.Range("Q6").Value = 3.793136154808667
$ws.Range("R6").Value = 34.13822539327801
$ws.Range("S6").Value = 0.001100207443118583
$ws.Range("T6").Value = 0.001100207443118583

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.028415333333333
$ws.Range("H7").Value = 3.085246
$ws.Range("I7").Value = 0.04565156193945813
$ws.Range("J7").Value = 0.04565156193945813
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.274860333333333
$ws.Range("N7").Value = 6.824581
$ws.Range("O7").Value = 0.01486427688856368
$ws.Range("P7").Value = 0.01486427688856368
$ws.Range("Q7").Value = 2.339501247991778
$ws.Range("R7").Value = 21.055511231926
$ws.Range("S7").Value = 0.0006785774570635208
$ws.Range("T7").Value = 0.0006785774570635208

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.925610333333334
$ws.Range("H8").Value = 14.776831
$ws.Range("I8").Value = 0.2186488259495045
$ws.Range("J8").Value = 0.2186488259495045
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 142.563596
$ws.Range("N8").Value = 427.690788
$ws.Range("O8").Value = 0.9315318106005318
$ws.Range("P8").Value = 0.9315318106005319
$ws.Range("Q8").Value = 702.2127216147587
$ws.Range("R8").Value = 6319.914494532829
$ws.Range("S8").Value = 0.2036783367224225
$ws.Range("T8").Value = 0.2036783367224225

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.925610333333334
$ws.Range("H9").Value = 14.776831
$ws.Range("I9").Value = 0.2186488259495045
$ws.Range("J9").Value = 0.2186488259495045
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.4502283333333333
$ws.Range("N9").Value = 1.350685
$ws.Range("O9").Value = 0.002941859116219682
$ws.Range("P9").Value = 0.002941859116219682
$ws.Range("Q9").Value = 2.217649331026111
$ws.Range("R9").Value = 19.958843979235
$ws.Range("S9").Value = 0.0006432340418702806
$ws.Range("T9").Value = 0.0006432340418702804

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.925610333333334
$ws.Range("H10").Value = 14.776831
$ws.Range("I10").Value = 0.2186488259495045
$ws.Range("J10").Value = 0.2186488259495045
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.398212333333333
$ws.Range("N10").Value = 4.194637
$ws.Range("O10").Value = 0.009136128037020016
$ws.Range("P10").Value = 0.009136128037020014
$ws.Range("Q10").Value = 6.887049117260779
$ws.Range("R10").Value = 61.98344205534701
$ws.Range("S10").Value = 0.001997603669018778
$ws.Range("T10").Value = 0.001997603669018778

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.925610333333334
$ws.Range("H11").Value = 14.776831
$ws.Range("I11").Value = 0.2186488259495045
$ws.Range("J11").Value = 0.2186488259495045
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.666884333333333
$ws.Range("N11").Value = 8.000653
$ws.Range("O11").Value = 0.01742582020512581
$ws.Range("P11").Value = 0.0174258202051258
$ws.Range("Q11").Value = 13.13603303007145
$ws.Range("R11").Value = 118.224297270643
$ws.Range("S11").Value = 0.003810135129057912
$ws.Range("T11").Value = 0.003810135129057911

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.925610333333334
$ws.Range("H12").Value = 14.776831
$ws.Range("I12").Value = 0.2186488259495045
$ws.Range("J12").Value = 0.2186488259495045
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.688331
$ws.Range("N12").Value = 11.064993
$ws.Range("O12").Value = 0.02410010515253888
$ws.Range("P12").Value = 0.02410010515253888
$ws.Range("Q12").Value = 18.16728128635367
$ws.Range("R12").Value = 163.505531577183
$ws.Range("S12").Value = 0.005269459696862231
$ws.Range("T12").Value = 0.00526945969686223

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.925610333333334
$ws.Range("H13").Value = 14.776831
$ws.Range("I13").Value = 0.2186488259495045
$ws.Range("J13").Value = 0.2186488259495045
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.274860333333333
$ws.Range("N13").Value = 6.824581
$ws.Range("O13").Value = 0.01486427688856368
$ws.Range("P13").Value = 0.01486427688856368
$ws.Range("Q13").Value = 11.20507556475678
$ws.Range("R13").Value = 100.845680082811
$ws.Range("S13").Value = 0.003250056690272803
$ws.Range("T13").Value = 0.003250056690272803

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.373401
$ws.Range("H14").Value = 13.120203
$ws.Range("I14").Value = 0.1941361434105301
$ws.Range("J14").Value = 0.1941361434105301
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 142.563596
$ws.Range("N14").Value = 427.690788
$ws.Range("O14").Value = 0.9315318106005318
$ws.Range("P14").Value = 0.9315318106005319
$ws.Range("Q14").Value = 623.487773309996
$ws.Range("R14").Value = 5611.389959789964
$ws.Range("S14").Value = 0.1808439931742156
$ws.Range("T14").Value = 0.1808439931742156

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.373401
$ws.Range("H15").Value = 13.120203
$ws.Range("I15").Value = 0.1941361434105301
$ws.Range("J15").Value = 0.1941361434105301
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.4502283333333333
$ws.Range("N15").Value = 1.350685
$ws.Range("O15").Value = 0.002941859116219682
$ws.Range("P15").Value = 0.002941859116219682
$ws.Range("Q15").Value = 1.969029043228333
$ws.Range("R15").Value = 17.721261389055
$ws.Range("S15").Value = 0.0005711211832799996
$ws.Range("T15").Value = 0.0005711211832799994

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.373401
$ws.Range("H16").Value = 13.120203
$ws.Range("I16").Value = 0.1941361434105301
$ws.Range("J16").Value = 0.1941361434105301
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.398212333333333
$ws.Range("N16").Value = 4.194637
$ws.Range("O16").Value = 0.009136128037020016
$ws.Range("P16").Value = 0.009136128037020014
$ws.Range("Q16").Value = 6.114943216812334
$ws.Range("R16").Value = 55.034488951311
$ws.Range("S16").Value = 0.001773652662811883
$ws.Range("T16").Value = 0.001773652662811882

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.373401
$ws.Range("H17").Value = 13.120203
$ws.Range("I17").Value = 0.1941361434105301
$ws.Range("J17").Value = 0.1941361434105301
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.666884333333333
$ws.Range("N17").Value = 8.000653
$ws.Range("O17").Value = 0.01742582020512581
$ws.Range("P17").Value = 0.0174258202051258
$ws.Range("Q17").Value = 11.66335461028434
$ws.Range("R17").Value = 104.970191492559
$ws.Range("S17").Value = 0.003382981530388417
$ws.Range("T17").Value = 0.003382981530388416

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 4.373401
$ws.Range("H18").Value = 13.120203
$ws.Range("I18").Value = 0.1941361434105301
$ws.Range("J18").Value = 0.1941361434105301
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 3.688331
$ws.Range("N18").Value = 11.064993
$ws.Range("O18").Value = 0.02410010515253888
$ws.Range("P18").Value = 0.02410010515253888
$ws.Range("Q18").Value = 16.130550483731
$ws.Range("R18").Value = 145.174954353579
$ws.Range("S18").Value = 0.004678701470102143
$ws.Range("T18").Value = 0.004678701470102142

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 4.373401
$ws.Range("H19").Value = 13.120203
$ws.Range("I19").Value = 0.1941361434105301
$ws.Range("J19").Value = 0.1941361434105301
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 2.274860333333333
$ws.Range("N19").Value = 6.824581
$ws.Range("O19").Value = 0.01486427688856368
$ws.Range("P19").Value = 0.01486427688856368
$ws.Range("Q19").Value = 9.948876456660333
$ws.Range("R19").Value = 89.539888109943
$ws.Range("S19").Value = 0.002885693389732027
$ws.Range("T19").Value = 0.002885693389732027

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 5.725097333333333
$ws.Range("H20").Value = 17.175292
$ws.Range("I20").Value = 0.2541382134735057
$ws.Range("J20").Value = 0.2541382134735057
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 142.563596
$ws.Range("N20").Value = 427.690788
$ws.Range("O20").Value = 0.9315318106005318
$ws.Range("P20").Value = 0.9315318106005319
$ws.Range("Q20").Value = 816.1904632900106
$ws.Range("R20").Value = 7345.714169610095
$ws.Range("S20").Value = 0.2367378301397593
$ws.Range("T20").Value = 0.2367378301397593

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 5.725097333333333
$ws.Range("H21").Value = 17.175292
$ws.Range("I21").Value = 0.2541382134735057
$ws.Range("J21").Value = 0.2541382134735057
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 0.4502283333333333
$ws.Range("N21").Value = 1.350685
$ws.Range("O21").Value = 0.002941859116219682
$ws.Range("P21").Value = 0.002941859116219682
$ws.Range("Q21").Value = 2.577601030557778
$ws.Range("R21").Value = 23.19840927502
$ws.Range("S21").Value = 0.0007476388200868165
$ws.Range("T21").Value = 0.0007476388200868164

# Row 22
$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 5.725097333333333
$ws.Range("H22").Value = 17.175292
$ws.Range("I22").Value = 0.2541382134735057
$ws.Range("J22").Value = 0.2541382134735057
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 1.398212333333333
$ws.Range("N22").Value = 4.194637
$ws.Range("O22").Value = 0.009136128037020016
$ws.Range("P22").Value = 0.009136128037020014
$ws.Range("Q22").Value = 8.004901701000446
$ws.Range("R22").Value = 72.044115309004
$ws.Range("S22").Value = 0.002321839257393474
$ws.Range("T22").Value = 0.002321839257393473

# Row 23
$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 5.725097333333333
$ws.Range("H23").Value = 17.175292
$ws.Range("I23").Value = 0.2541382134735057
$ws.Range("J23").Value = 0.2541382134735057
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 2.666884333333333
$ws.Range("N23").Value = 8.000653
$ws.Range("O23").Value = 0.01742582020512581
$ws.Range("P23").Value = 0.0174258202051258
$ws.Range("Q23").Value = 15.26817238507511
$ws.Range("R23").Value = 137.413551465676
$ws.Range("S23").Value = 0.004428566815241191
$ws.Range("T23").Value = 0.004428566815241191

# Row 24
$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 5.725097333333333
$ws.Range("H24").Value = 17.175292
$ws.Range("I24").Value = 0.2541382134735057
$ws.Range("J24").Value = 0.2541382134735057
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 3.688331
$ws.Range("N24").Value = 11.064993
$ws.Range("O24").Value = 0.02410010515253888
$ws.Range("P24").Value = 0.02410010515253888
$ws.Range("Q24").Value = 21.11605397255067
$ws.Range("R24").Value = 190.044485752956
$ws.Range("S24").Value = 0.006124757667989861
$ws.Range("T24").Value = 0.006124757667989861

# Row 25
$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 5.725097333333333
$ws.Range("H25").Value = 17.175292
$ws.Range("I25").Value = 0.2541382134735057
$ws.Range("J25").Value = 0.2541382134735057
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 2.274860333333333
$ws.Range("N25").Value = 6.824581
$ws.Range("O25").Value = 0.01486427688856368
$ws.Range("P25").Value = 0.01486427688856368
$ws.Range("Q25").Value = 13.02379682807244
$ws.Range("R25").Value = 117.214171452652
$ws.Range("S25").Value = 0.003777580773035094
$ws.Range("T25").Value = 0.003777580773035094

# Row 26
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 1.796371666666667
$ws.Range("H26").Value = 5.389115
$ws.Range("I26").Value = 0.07974129687595831
$ws.Range("J26").Value = 0.07974129687595832
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 142.563596
$ws.Range("N26").Value = 427.690788
$ws.Range("O26").Value = 0.9315318106005318
$ws.Range("P26").Value = 0.9315318106005319
$ws.Range("Q26").Value = 256.0972045525133
$ws.Range("R26").Value = 2304.87484097262
$ws.Range("S26").Value = 0.07428155465849598
$ws.Range("T26").Value = 0.074281554658496

# Row 27
$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 1.796371666666667
$ws.Range("H27").Value = 5.389115
$ws.Range("I27").Value = 0.07974129687595831
$ws.Range("J27").Value = 0.07974129687595832
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 0.4502283333333333
$ws.Range("N27").Value = 1.350685
$ws.Range("O27").Value = 0.002941859116219682
$ws.Range("P27").Value = 0.002941859116219682
$ws.Range("Q27").Value = 0.8087774215305555
$ws.Range("R27").Value = 7.278996793775
$ws.Range("S27").Value = 0.000234587661153718
$ws.Range("T27").Value = 0.000234587661153718

# Row 28
$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 1.796371666666667
$ws.Range("H28").Value = 5.389115
$ws.Range("I28").Value = 0.07974129687595831
$ws.Range("J28").Value = 0.07974129687595832
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 1.398212333333333
$ws.Range("N28").Value = 4.194637
$ws.Range("O28").Value = 0.009136128037020016
$ws.Range("P28").Value = 0.009136128037020014
$ws.Range("Q28").Value = 2.511709019583889
$ws.Range("R28").Value = 22.605381176255
$ws.Range("S28").Value = 0.0007285266980967793
$ws.Range("T28").Value = 0.0007285266980967793

# Row 29
$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 1.796371666666667
$ws.Range("H29").Value = 5.389115
$ws.Range("I29").Value = 0.07974129687595831
$ws.Range("J29").Value = 0.07974129687595832
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 2.666884333333333
$ws.Range("N29").Value = 8.000653
$ws.Range("O29").Value = 0.01742582020512581
$ws.Range("P29").Value = 0.0174258202051258
$ws.Range("Q29").Value = 4.790715454677223
$ws.Range("R29").Value = 43.116439092095
$ws.Range("S29").Value = 0.00138955750228401
$ws.Range("T29").Value = 0.00138955750228401

# Row 30
$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 1.796371666666667
$ws.Range("H30").Value = 5.389115
$ws.Range("I30").Value = 0.07974129687595831
$ws.Range("J30").Value = 0.07974129687595832
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 3.688331
$ws.Range("N30").Value = 11.064993
$ws.Range("O30").Value = 0.02410010515253888
$ws.Range("P30").Value = 0.02410010515253888
$ws.Range("Q30").Value = 6.625613305688334
$ws.Range("R30").Value = 59.63051975119501
$ws.Range("S30").Value = 0.001921773639710415
$ws.Range("T30").Value = 0.001921773639710415

# Row 31
$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 1.796371666666667
$ws.Range("H31").Value = 5.389115
$ws.Range("I31").Value = 0.07974129687595831
$ws.Range("J31").Value = 0.07974129687595832
$ws.Range("K31").Value = 3
$ws.Range("M31").Value = 2.274860333333333
$ws.Range("N31").Value = 6.824581
$ws.Range("O31").Value = 0.01486427688856368
$ws.Range("P31").Value = 0.01486427688856368
$ws.Range("Q31").Value = 4.086494648423889
$ws.Range("R31").Value = 36.77845183581501
$ws.Range("S31").Value = 0.001184019639955765
$ws.Range("T31").Value = 0.001184019639955765

# Row 32
$ws.Range("E32").Value = 3
$ws.Range("G32").Value = 4.678599333333334
$ws.Range("H32").Value = 14.035798
$ws.Range("I32").Value = 0.2076839583510432
$ws.Range("J32").Value = 0.2076839583510431
$ws.Range("K32").Value = 3
$ws.Range("M32").Value = 142.563596
$ws.Range("N32").Value = 427.690788
$ws.Range("O32").Value = 0.9315318106005318
$ws.Range("P32").Value = 0.9315318106005319
$ws.Range("Q32").Value = 666.9979452032027
$ws.Range("R32").Value = 6002.981506828824
$ws.Range("S32").Value = 0.1934642137554327
$ws.Range("T32").Value = 0.1934642137554327

# Row 33
$ws.Range("E33").Value = 3
$ws.Range("G33").Value = 4.678599333333334
$ws.Range("H33").Value = 14.035798
$ws.Range("I33").Value = 0.2076839583510432
$ws.Range("J33").Value = 0.2076839583510431
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 0.4502283333333333
$ws.Range("N33").Value = 1.350685
$ws.Range("O33").Value = 0.002941859116219682
$ws.Range("P33").Value = 0.002941859116219682
$ws.Range("Q33").Value = 2.106437980181111
$ws.Range("R33").Value = 18.95794182163
$ws.Range("S33").Value = 0.0006109769461676051
$ws.Range("T33").Value = 0.0006109769461676049

# Row 34
$ws.Range("E34").Value = 3
$ws.Range("G34").Value = 4.678599333333334
$ws.Range("H34").Value = 14.035798
$ws.Range("I34").Value = 0.2076839583510432
$ws.Range("J34").Value = 0.2076839583510431
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 1.398212333333333
$ws.Range("N34").Value = 4.194637
$ws.Range("O34").Value = 0.009136128037020016
$ws.Range("P34").Value = 0.009136128037020014
$ws.Range("Q34").Value = 6.541675290591779
$ws.Range("R34").Value = 58.875077615326
$ws.Range("S34").Value = 0.001897427234730263
$ws.Range("T34").Value = 0.001897427234730262

# Row 35
$ws.Range("E35").Value = 3
$ws.Range("G35").Value = 4.678599333333334
$ws.Range("H35").Value = 14.035798
$ws.Range("I35").Value = 0.2076839583510432
$ws.Range("J35").Value = 0.2076839583510431
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 2.666884333333333
$ws.Range("N35").Value = 8.000653
$ws.Range("O35").Value = 0.01742582020512581
$ws.Range("P35").Value = 0.0174258202051258
$ws.Range("Q35").Value = 12.47728326401045
$ws.Range("R35").Value = 112.295549376094
$ws.Range("S35").Value = 0.003619063317714115
$ws.Range("T35").Value = 0.003619063317714113

# Row 36
$ws.Range("E36").Value = 3
$ws.Range("G36").Value = 4.678599333333334
$ws.Range("H36").Value = 14.035798
$ws.Range("I36").Value = 0.2076839583510432
$ws.Range("J36").Value = 0.2076839583510431
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 3.688331
$ws.Range("N36").Value = 11.064993
$ws.Range("O36").Value = 0.02410010515253888
$ws.Range("P36").Value = 0.02410010515253888
$ws.Range("Q36").Value = 17.25622295771267
$ws.Range("R36").Value = 155.306006619414
$ws.Range("S36").Value = 0.005005205234755645
$ws.Range("T36").Value = 0.005005205234755644

# Row 37
$ws.Range("E37").Value = 3
$ws.Range("G37").Value = 4.678599333333334
$ws.Range("H37").Value = 14.035798
$ws.Range("I37").Value = 0.2076839583510432
$ws.Range("J37").Value = 0.2076839583510431
$ws.Range("K37").Value = 3
$ws.Range("M37").Value = 2.274860333333333
$ws.Range("N37").Value = 6.824581
$ws.Range("O37").Value = 0.01486427688856368
$ws.Range("P37").Value = 0.01486427688856368
$ws.Range("Q37").Value = 10.64316003895978
$ws.Range("R37").Value = 95.788440350638
$ws.Range("S37").Value = 0.003087071862242833
$ws.Range("T37").Value = 0.003087071862242833
